$wb = $excel.ActiveWorkbook

# Rename the "rooms" sheet to "beds" (terminology change: room -> bed)
$ws = $wb.Worksheets.Item("rooms")
$ws.Name = "beds"

# Update header row terminology on the renamed sheet
$ws.Range("A1").Value = "all_beds"
$ws.Range("B1").Value = "new_beds"
$ws.Range("C1").Value = "old_beds"
$ws.Range("E1").Value = "new_beds_service"
$ws.Range("F1").Value = "old_beds_service"
$ws.Range("G1").Value = "beds_capacities"

# Make the "beds" sheet the active tab (was previously "babies")
$ws.Activate()
